$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pool")

$ws.Range("AP1").Value = "VL"
$ws.Range("AQ1").Value = "SL"
$ws.Range("AP2").Value = "15"
$ws.Range("AQ2").Value = "15"

$ws.Range("AP1:AQ2").Select()
